$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.336433666666667
$ws.Range("H2").Value = 16.009301
$ws.Range("I2").Value = 0.09554123656860429
$ws.Range("J2").Value = 0.09601942232049432
$ws.Range("M2").Value = 3.058122
$ws.Range("N2").Value = 9.174365999999999
$ws.Range("O2").Value = 0.1133977796540004
$ws.Range("P2").Value = 0.1165634753455787
$ws.Range("Q2").Value = 16.319465197574
$ws.Range("R2").Value = 146.875186778166
$ws.Range("S2").Value = 0.01083416409227731
$ws.Range("T2").Value = 0.01119235756635165
$ws.Range("G3").Value = 5.336433666666667
$ws.Range("H3").Value = 16.009301
$ws.Range("I3").Value = 0.09554123656860429
$ws.Range("J3").Value = 0.09601942232049432
$ws.Range("O3").Value = 0.4569904442326647
$ws.Range("P3").Value = 0.4697481250692204
$ws.Range("Q3").Value = 65.76706945263157
$ws.Range("R3").Value = 591.903625073684
$ws.Range("S3").Value = 0.04366143214202459
$ws.Range("T3").Value = 0.04510494360528185
$ws.Range("G4").Value = 5.336433666666667
$ws.Range("H4").Value = 16.009301
$ws.Range("I4").Value = 0.09554123656860429
$ws.Range("J4").Value = 0.09601942232049432
$ws.Range("M4").Value = 6.332890333333334
$ws.Range("N4").Value = 18.998671
$ws.Range("O4").Value = 0.2348289906656054
$ws.Range("P4").Value = 0.241384649217969
$ws.Range("Q4").Value = 33.7950491821079
$ws.Range("R4").Value = 304.1554426389711
$ws.Range("S4").Value = 0.02243585215034917
$ws.Range("T4").Value = 0.02317761457494454
$ws.Range("G5").Value = 5.336433666666667
$ws.Range("H5").Value = 16.009301
$ws.Range("I5").Value = 0.09554123656860429
$ws.Range("J5").Value = 0.09601942232049432
$ws.Range("M5").Value = 2.1972435
$ws.Range("N5").Value = 4.394487
$ws.Range("O5").Value = 0.08147566848516331
$ws.Range("P5").Value = 0.05583346871936067
$ws.Range("Q5").Value = 11.7254441872645
$ws.Range("R5").Value = 70.352665123587
$ws.Range("S5").Value = 0.007784286117326165
$ws.Range("T5").Value = 0.005361097412582401
$ws.Range("G6").Value = 5.336433666666667
$ws.Range("H6").Value = 16.009301
$ws.Range("I6").Value = 0.09554123656860429
$ws.Range("J6").Value = 0.09601942232049432
$ws.Range("M6").Value = 3.055677
$ws.Range("N6").Value = 9.167031
$ws.Range("O6").Value = 0.1133071169625662
$ws.Range("P6").Value = 0.1164702816478714
$ws.Range("Q6").Value = 16.306417617259
$ws.Range("R6").Value = 146.757758555331
$ws.Range("S6").Value = 0.01082550206662706
$ws.Range("T6").Value = 0.01118340916133388
$ws.Range("I7").Value = 0.0653459693259494
$ws.Range("J7").Value = 0.06567302717654226
$ws.Range("M7").Value = 3.058122
$ws.Range("N7").Value = 9.174365999999999
$ws.Range("O7").Value = 0.1133977796540004
$ws.Range("P7").Value = 0.1165634753455787
$ws.Range("Q7").Value = 11.161790557848
$ws.Range("R7").Value = 100.456115020632
$ws.Range("S7").Value = 0.007410087830901077
$ws.Range("T7").Value = 0.007655076284162402
$ws.Range("I8").Value = 0.0653459693259494
$ws.Range("J8").Value = 0.06567302717654226
$ws.Range("O8").Value = 0.4569904442326647
$ws.Range("P8").Value = 0.4697481250692204
$ws.Range("S8").Value = 0.0298624835510797
$ws.Range("T8").Value = 0.03084978138380068
$ws.Range("I9").Value = 0.0653459693259494
$ws.Range("J9").Value = 0.06567302717654226
$ws.Range("M9").Value = 6.332890333333334
$ws.Range("N9").Value = 18.998671
$ws.Range("O9").Value = 0.2348289906656054
$ws.Range("P9").Value = 0.241384649217969
$ws.Range("Q9").Value = 23.11431510138801
$ws.Range("R9").Value = 208.028835912492
$ws.Range("S9").Value = 0.01534512802087831
$ws.Range("T9").Value = 0.0158524606280918
$ws.Range("I10").Value = 0.0653459693259494
$ws.Range("J10").Value = 0.06567302717654226
$ws.Range("M10").Value = 2.1972435
$ws.Range("N10").Value = 4.394487
$ws.Range("O10").Value = 0.08147566848516331
$ws.Range("P10").Value = 0.05583346871936067
$ws.Range("Q10").Value = 8.019683894753999
$ws.Range("R10").Value = 48.118103368524
$ws.Range("S10").Value = 0.005324106533642704
$ws.Range("T10").Value = 0.003666752908567195
$ws.Range("I11").Value = 0.0653459693259494
$ws.Range("J11").Value = 0.06567302717654226
$ws.Range("M11").Value = 3.055677
$ws.Range("N11").Value = 9.167031
$ws.Range("O11").Value = 0.1133071169625662
$ws.Range("P11").Value = 0.1164702816478714
$ws.Range("Q11").Value = 11.152866591468
$ws.Range("R11").Value = 100.375799323212
$ws.Range("S11").Value = 0.007404163389447613
$ws.Range("T11").Value = 0.007648955971920191
$ws.Range("G12").Value = 22.33109633333333
$ws.Range("H12").Value = 66.993289
$ws.Range("I12").Value = 0.3998064420712607
$ws.Range("J12").Value = 0.4018074811092581
$ws.Range("M12").Value = 3.058122
$ws.Range("N12").Value = 9.174365999999999
$ws.Range("O12").Value = 0.1133977796540004
$ws.Range("P12").Value = 0.1165634753455787
$ws.Range("Q12").Value = 68.291216981086
$ws.Range("R12").Value = 614.620952829774
$ws.Range("S12").Value = 0.04533716282224669
$ws.Range("T12").Value = 0.04683607641794807
$ws.Range("G13").Value = 22.33109633333333
$ws.Range("H13").Value = 66.993289
$ws.Range("I13").Value = 0.3998064420712607
$ws.Range("J13").Value = 0.4018074811092581
$ws.Range("O13").Value = 0.4569904442326647
$ws.Range("P13").Value = 0.4697481250692204
$ws.Range("Q13").Value = 275.2120339622085
$ws.Range("R13").Value = 2476.908305659876
$ws.Range("S13").Value = 0.1827077235692266
$ws.Range("T13").Value = 0.1887483108898602
$ws.Range("G14").Value = 22.33109633333333
$ws.Range("H14").Value = 66.993289
$ws.Range("I14").Value = 0.3998064420712607
$ws.Range("J14").Value = 0.4018074811092581
$ws.Range("M14").Value = 6.332890333333334
$ws.Range("N14").Value = 18.998671
$ws.Range("O14").Value = 0.2348289906656054
$ws.Range("P14").Value = 0.241384649217969
$ws.Range("Q14").Value = 141.4203841021021
$ws.Range("R14").Value = 1272.783456918919
$ws.Range("S14").Value = 0.09388614325320098
$ws.Range("T14").Value = 0.09699015788071394
$ws.Range("G15").Value = 22.33109633333333
$ws.Range("H15").Value = 66.993289
$ws.Range("I15").Value = 0.3998064420712607
$ws.Range("J15").Value = 0.4018074811092581
$ws.Range("M15").Value = 2.1972435
$ws.Range("N15").Value = 4.394487
$ws.Range("O15").Value = 0.08147566848516331
$ws.Range("P15").Value = 0.05583346871936067
$ws.Range("Q15").Value = 49.0668562662905
$ws.Range("R15").Value = 294.401137597743
$ws.Range("S15").Value = 0.03257449713243069
$ws.Range("T15").Value = 0.02243430542771886
$ws.Range("G16").Value = 22.33109633333333
$ws.Range("H16").Value = 66.993289
$ws.Range("I16").Value = 0.3998064420712607
$ws.Range("J16").Value = 0.4018074811092581
$ws.Range("M16").Value = 3.055677
$ws.Range("N16").Value = 9.167031
$ws.Range("O16").Value = 0.1133071169625662
$ws.Range("P16").Value = 0.1164702816478714
$ws.Range("Q16").Value = 68.23661745055099
$ws.Range("R16").Value = 614.129557054959
$ws.Range("S16").Value = 0.04530091529415579
$ws.Range("T16").Value = 0.04679863049301706
$ws.Range("G17").Value = 0.8344860000000001
$ws.Range("H17").Value = 1.668972
$ws.Range("I17").Value = 0.01494028209086487
$ws.Range("J17").Value = 0.01001003899602363
$ws.Range("M17").Value = 3.058122
$ws.Range("N17").Value = 9.174365999999999
$ws.Range("O17").Value = 0.1133977796540004
$ws.Range("P17").Value = 0.1165634753455787
$ws.Range("Q17").Value = 2.551959995292
$ws.Range("R17").Value = 15.311759971752
$ws.Range("S17").Value = 0.001694194816508502
$ws.Range("T17").Value = 0.001166804933721281
$ws.Range("G18").Value = 0.8344860000000001
$ws.Range("H18").Value = 1.668972
$ws.Range("I18").Value = 0.01494028209086487
$ws.Range("J18").Value = 0.01001003899602363
$ws.Range("O18").Value = 0.4569904442326647
$ws.Range("P18").Value = 0.4697481250692204
$ws.Range("Q18").Value = 10.284340094408
$ws.Range("R18").Value = 61.70604056644801
$ws.Range("S18").Value = 0.006827566149665659
$ws.Range("T18").Value = 0.00470219705025188
$ws.Range("G19").Value = 0.8344860000000001
$ws.Range("H19").Value = 1.668972
$ws.Range("I19").Value = 0.01494028209086487
$ws.Range("J19").Value = 0.01001003899602363
$ws.Range("M19").Value = 6.332890333333334
$ws.Range("N19").Value = 18.998671
$ws.Range("O19").Value = 0.2348289906656054
$ws.Range("P19").Value = 0.241384649217969
$ws.Range("Q19").Value = 5.284708322702001
$ws.Range("R19").Value = 31.70824993621201
$ws.Range("S19").Value = 0.003508411363657217
$ws.Range("T19").Value = 0.002416269751713353
$ws.Range("G20").Value = 0.8344860000000001
$ws.Range("H20").Value = 1.668972
$ws.Range("I20").Value = 0.01494028209086487
$ws.Range("J20").Value = 0.01001003899602363
$ws.Range("M20").Value = 2.1972435
$ws.Range("N20").Value = 4.394487
$ws.Range("O20").Value = 0.08147566848516331
$ws.Range("P20").Value = 0.05583346871936067
$ws.Range("Q20").Value = 1.833568939341
$ws.Range("R20").Value = 7.334275757364
$ws.Range("S20").Value = 0.001217269470710128
$ws.Range("T20").Value = 0.0005588951991640657
$ws.Range("G21").Value = 0.8344860000000001
$ws.Range("H21").Value = 1.668972
$ws.Range("I21").Value = 0.01494028209086487
$ws.Range("J21").Value = 0.01001003899602363
$ws.Range("M21").Value = 3.055677
$ws.Range("N21").Value = 9.167031
$ws.Range("O21").Value = 0.1133071169625662
$ws.Range("P21").Value = 0.1164702816478714
$ws.Range("Q21").Value = 2.549919677022
$ws.Range("R21").Value = 15.299518062132
$ws.Range("S21").Value = 0.001692840290323359
$ws.Range("T21").Value = 0.001165872061173048
$ws.Range("G22").Value = 23.70286866666666
$ws.Range("H22").Value = 71.10860599999999
$ws.Range("I22").Value = 0.4243660699433207
$ws.Range("J22").Value = 0.4264900303976816
$ws.Range("M22").Value = 3.058122
$ws.Range("N22").Value = 9.174365999999999
$ws.Range("O22").Value = 0.1133977796540004
$ws.Range("P22").Value = 0.1165634753455787
$ws.Range("Q22").Value = 72.48626413264398
$ws.Range("R22").Value = 652.3763771937959
$ws.Range("S22").Value = 0.04812217009206678
$ws.Range("T22").Value = 0.04971316014339526
$ws.Range("G23").Value = 23.70286866666666
$ws.Range("H23").Value = 71.10860599999999
$ws.Range("I23").Value = 0.4243660699433207
$ws.Range("J23").Value = 0.4264900303976816
$ws.Range("O23").Value = 0.4569904442326647
$ws.Range("P23").Value = 0.4697481250692204
$ws.Range("Q23").Value = 292.1179775108116
$ws.Range("R23").Value = 2629.061797597304
$ws.Range("S23").Value = 0.1939312388206682
$ws.Range("T23").Value = 0.2003428921400257
$ws.Range("G24").Value = 23.70286866666666
$ws.Range("H24").Value = 71.10860599999999
$ws.Range("I24").Value = 0.4243660699433207
$ws.Range("J24").Value = 0.4264900303976816
$ws.Range("M24").Value = 6.332890333333334
$ws.Range("N24").Value = 18.998671
$ws.Range("O24").Value = 0.2348289906656054
$ws.Range("P24").Value = 0.241384649217969
$ws.Range("Q24").Value = 150.1076678514029
$ws.Range("R24").Value = 1350.969010662626
$ws.Range("S24").Value = 0.0996534558775197
$ws.Range("T24").Value = 0.1029481463825053
$ws.Range("G25").Value = 23.70286866666666
$ws.Range("H25").Value = 71.10860599999999
$ws.Range("I25").Value = 0.4243660699433207
$ws.Range("J25").Value = 0.4264900303976816
$ws.Range("M25").Value = 2.1972435
$ws.Range("N25").Value = 4.394487
$ws.Range("O25").Value = 0.08147566848516331
$ws.Range("P25").Value = 0.05583346871936067
$ws.Range("Q25").Value = 52.08097410918699
$ws.Range("R25").Value = 312.485844655122
$ws.Range("S25").Value = 0.03457550923105362
$ws.Range("T25").Value = 0.02381241777132814
$ws.Range("G26").Value = 23.70286866666666
$ws.Range("H26").Value = 71.10860599999999
$ws.Range("I26").Value = 0.4243660699433207
$ws.Range("J26").Value = 0.4264900303976816
$ws.Range("M26").Value = 3.055677
$ws.Range("N26").Value = 9.167031
$ws.Range("O26").Value = 0.1133071169625662
$ws.Range("P26").Value = 0.1164702816478714
$ws.Range("Q26").Value = 72.42831061875398
$ws.Range("R26").Value = 651.8547955687859
$ws.Range("S26").Value = 0.04808369592201239
$ws.Range("T26").Value = 0.04967341396042722
